$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "Removed the guesses for the B spectral types from the tables"
# Rows 38-43 held guessed parameters for B0I, B1I, B0III, B1III, B0V, B1V.
# B0I/B1I (38/39) and B0V/B1V (42/43) are dropped entirely; B0III/B1III
# (40/41) keep their row but are wiped down to a single empty, number-
# formatted B cell (the row/style survives, the data does not).

# Row 38 (B0I) and row 39 (B1I): clear entirely
$ws.Range("A38:L39").ClearContents() | Out-Null

# Row 40 (B0III): keep only the styled, now-empty B cell
$ws.Range("A40").ClearContents() | Out-Null
$ws.Range("C40:L40").ClearContents() | Out-Null
$ws.Range("B40").ClearContents() | Out-Null

# Row 41 (B1III): keep only the styled, now-empty B cell
$ws.Range("A41").ClearContents() | Out-Null
$ws.Range("C41:L41").ClearContents() | Out-Null
$ws.Range("B41").ClearContents() | Out-Null

# Row 42 (B0V) and row 43 (B1V): clear entirely
$ws.Range("A42:L43").ClearContents() | Out-Null

# Update the visible selection to match the saved state after the edit
$ws.Range("J46").Select() | Out-Null
